$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New column D values (historical_growth_revenue_last_5_years) for data rows
$ws.Range("D2").Value = -0.0372
$ws.Range("D3").Value = -0.0372

foreach ($r in 2,3) {
    $ws.Range("G$r").Value = -0.8216867469879519
    $ws.Range("H$r").Value = -0.8216867469879519
    $ws.Range("I$r").Value = -0.8313253012048193
    $ws.Range("J$r").Value = -0.8313253012048193
    $ws.Range("K$r").Value = -32.8
    $ws.Range("L$r").Value = -0.7903614457831325
    $ws.Range("M$r").Value = 0.112
    $ws.Range("N$r").Value = 0.007777777777777778
    $ws.Range("O$r").Value = -0.003414634146341464
    $ws.Range("S$r").Value = 0.112
    $ws.Range("U$r").Value = 8.6
    $ws.Range("V$r").Value = 0.5972222222222222
    $ws.Range("W$r").Value = -0.6295585412667946
    $ws.Range("X$r").Value = 0.0697865179974401
    $ws.Range("Y$r").Value = -0.6993450592642346
    $ws.Range("Z$r").Value = 1.039579158316633
    $ws.Range("AA$r").Value = -0.8642284569138277
    $ws.Range("AB$r").Value = 0.06347585515073567
    $ws.Range("AC$r").Value = -0.9277043120645634
    $ws.Range("AD$r").Value = 3.99
    $ws.Range("AF$r").Value = 3.99
    $ws.Range("AG$r").Value = -4.609999999999999
    $ws.Range("AH$r").Value = 0.2169657422512235
    $ws.Range("AI$r").Value = 0.2831795599716111
    $ws.Range("AJ$r").Value = -0.4708886618998978
    $ws.Range("AK$r").Value = -0.8397085610200363
    $ws.Range("AL$r").Value = 0.48
    $ws.Range("AM$r").Value = 0.48
    $ws.Range("AN$r").Value = -0.1170087976539589
    $ws.Range("AO$r").Value = -71.875
    $ws.Range("AP$r").Value = 0.1351906158357771
    $ws.Range("AQ$r").Value = -71.875
}
